# Rename the sheet from "Validated Addresses" to "Validated_Results"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Validated_Results"

# Write the header row
$ws.Range("A1").Value = "address"
$ws.Range("B1").Value = "validatedAddresses"
$ws.Range("C1").Value = "coordinates"
$ws.Range("D1").Value = "resolutionQuality"
$ws.Range("E1").Value = "taxAuthorities"

# Write the two data rows (only column D has data)
$ws.Range("D2").Value = "Intersection"
$ws.Range("D3").Value = "Intersection"
